# Update Leave Card 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------
# Step 1: Insert two new rows into the leave table.
#   - one new row right after row 500 (Oct 2023) for an extra
#     Sick-Leave entry dated 10/16/2023
#   - one new row that becomes the "2024" year-separator row
# After each insert we immediately resize Table1 to include the
# new row and restore the structured-reference formula on what is
# now the table's last row (Excel rewrites it into a plain,
# non-table formula - and leaves it outside the table range -
# when a row is inserted above the table's final row).
# ---------------------------------------------------------------
$ws.Rows.Item(501).Insert()
$lo.Resize($ws.Range("A8:K603"))
$ws.Range("G603").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$ws.Rows.Item(504).Insert()
$lo.Resize($ws.Range("A8:K604"))
$ws.Range("G604").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------
# Step 2: Re-apply the correct cell formatting to the freshly
# inserted rows (Excel gives newly inserted rows generic/new
# styles); copy it from a clean, untouched template row that
# already has the right per-column formatting.
# ---------------------------------------------------------------
$ws.Range("A510:K510").Copy()
$ws.Range("A501:K501").PasteSpecial(-4122)
$ws.Range("A504:K504").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the calculated-column formula (structured references) on
# the newly inserted rows, since paste-formats does not bring the
# formula back.
$ws.Range("G501").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G504").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------
# Step 3: Populate row 500 (October 2023) with a new SL entry
# ---------------------------------------------------------------
$ws.Range("B500").Value = "SL(1-0-0)"
$ws.Range("C500").Value = 1.25
$ws.Range("H500").Value = 1
$ws.Range("K499").Copy()
$ws.Range("K500").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K500").Value = 45208

# ---------------------------------------------------------------
# Step 4: Populate the newly inserted row 501 (no date, SL entry)
# ---------------------------------------------------------------
$ws.Range("B501").Value = "SL(1-0-0)"
$ws.Range("H501").Value = 1
$ws.Range("K499").Copy()
$ws.Range("K501").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K501").Value = 45215

# ---------------------------------------------------------------
# Step 5: Row 502 (November 2023, shifted from old row 501) gets
# an EARNED value
# ---------------------------------------------------------------
$ws.Range("C502").Value = 1.25

# ---------------------------------------------------------------
# Step 6: Row 503 (December 2023, shifted from old row 502) gets
# an SL entry
# ---------------------------------------------------------------
$ws.Range("B503").Value = "SL(1-0-0)"
$ws.Range("H503").Value = 1

# ---------------------------------------------------------------
# Step 7: Row 504 becomes the "2024" year separator row
# ---------------------------------------------------------------
$ws.Range("A10").Copy()
$ws.Range("A504").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A504").Value = "'2024"

# ---------------------------------------------------------------
# Step 8: Recalculate so dependent formulas (BALANCE columns, etc.)
# pick up the new figures
# ---------------------------------------------------------------
$excel.CalculateFullRebuild()

# Restore the selection to reflect where the user last clicked
$ws.Range("I503").Select()
